$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.670167666666667
$ws.Range("H2").Value = 20.010503
$ws.Range("I2").Value = 0.0423069620011633
$ws.Range("J2").Value = 0.0423069620011633
$ws.Range("M2").Value = 22.618885
$ws.Range("N2").Value = 67.856655
$ws.Range("O2").Value = 0.9027998993061069
$ws.Range("P2").Value = 0.902799899306107
$ws.Range("Q2").Value = 150.8717553830517
$ws.Range("R2").Value = 1357.845798447465
$ws.Range("S2").Value = 0.03819472103459751
$ws.Range("T2").Value = 0.03819472103459752
$ws.Range("G3").Value = 6.670167666666667
$ws.Range("H3").Value = 20.010503
$ws.Range("I3").Value = 0.0423069620011633
$ws.Range("J3").Value = 0.0423069620011633
$ws.Range("O3").Value = 0.08600330007856447
$ws.Range("P3").Value = 0.08600330007856449
$ws.Range("Q3").Value = 14.37247485468411
$ws.Range("R3").Value = 129.352273692157
$ws.Range("S3").Value = 0.003638538348398471
$ws.Range("T3").Value = 0.003638538348398472
$ws.Range("G4").Value = 6.670167666666667
$ws.Range("H4").Value = 20.010503
$ws.Range("I4").Value = 0.0423069620011633
$ws.Range("J4").Value = 0.0423069620011633
$ws.Range("M4").Value = 0.2805263333333333
$ws.Range("N4").Value = 0.841579
$ws.Range("O4").Value = 0.0111968006153285
$ws.Range("P4").Value = 0.01119680061532851
$ws.Range("Q4").Value = 1.871157678248556
$ws.Range("R4").Value = 16.840419104237
$ws.Range("S4").Value = 0.0004737026181673049
$ws.Range("T4").Value = 0.0004737026181673049
$ws.Range("I5").Value = 0.9513278459982415
$ws.Range("J5").Value = 0.9513278459982416
$ws.Range("M5").Value = 22.618885
$ws.Range("N5").Value = 67.856655
$ws.Range("O5").Value = 0.9027998993061069
$ws.Range("P5").Value = 0.902799899306107
$ws.Range("Q5").Value = 3392.550428617059
$ws.Range("R5").Value = 30532.95385755353
$ws.Range("S5").Value = 0.8588586835743081
$ws.Range("T5").Value = 0.8588586835743083
$ws.Range("I6").Value = 0.9513278459982415
$ws.Range("J6").Value = 0.9513278459982416
$ws.Range("O6").Value = 0.08600330007856447
$ws.Range("P6").Value = 0.08600330007856449
$ws.Range("S6").Value = 0.08181733421248112
$ws.Range("T6").Value = 0.08181733421248116
$ws.Range("I7").Value = 0.9513278459982415
$ws.Range("J7").Value = 0.9513278459982416
$ws.Range("M7").Value = 0.2805263333333333
$ws.Range("N7").Value = 0.841579
$ws.Range("O7").Value = 0.0111968006153285
$ws.Range("P7").Value = 0.01119680061532851
$ws.Range("Q7").Value = 42.07544856381612
$ws.Range("R7").Value = 378.679037074345
$ws.Range("S7").Value = 0.01065182821145225
$ws.Range("T7").Value = 0.01065182821145225
$ws.Range("G8").Value = 1.003544
$ws.Range("H8").Value = 3.010632
$ws.Range("I8").Value = 0.0063651920005952
$ws.Range("J8").Value = 0.006365192000595201
$ws.Range("M8").Value = 22.618885
$ws.Range("N8").Value = 67.856655
$ws.Range("O8").Value = 0.9027998993061069
$ws.Range("P8").Value = 0.902799899306107
$ws.Range("Q8").Value = 22.69904632844
$ws.Range("R8").Value = 204.29141695596
$ws.Range("S8").Value = 0.005746494697201384
$ws.Range("T8").Value = 0.005746494697201385
$ws.Range("G9").Value = 1.003544
$ws.Range("H9").Value = 3.010632
$ws.Range("I9").Value = 0.0063651920005952
$ws.Range("J9").Value = 0.006365192000595201
$ws.Range("O9").Value = 0.08600330007856447
$ws.Range("P9").Value = 0.08600330007856449
$ws.Range("Q9").Value = 2.162376064045333
$ws.Range("R9").Value = 19.461384576408
$ws.Range("S9").Value = 0.000547427517684867
$ws.Range("T9").Value = 0.0005474275176848673
$ws.Range("G10").Value = 1.003544
$ws.Range("H10").Value = 3.010632
$ws.Range("I10").Value = 0.0063651920005952
$ws.Range("J10").Value = 0.006365192000595201
$ws.Range("M10").Value = 0.2805263333333333
$ws.Range("N10").Value = 0.841579
$ws.Range("O10").Value = 0.0111968006153285
$ws.Range("P10").Value = 0.01119680061532851
$ws.Range("Q10").Value = 0.2815205186586667
$ws.Range("R10").Value = 2.533684667928
$ws.Range("S10").Value = 0.00007126978570894841
$ws.Range("T10").Value = 0.00007126978570894844
